$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 2 (shifts existing rows 2 and 3 down to 3 and 4)
$ws.Rows.Item(2).Insert()

# Copy the date-number formatting (style) from the date cell now in row 3
# onto the new row 2 date cell, so it reuses the same style index.
$ws.Range("C3").Copy($ws.Range("C2"))

# Populate the newly inserted row 2 with the "sex" entry
$ws.Range("A2").Value = "sex"
$ws.Range("B2").Value = 30
$ws.Range("C2").Value = 45917.22928240741

# Update the (now shifted down) former "Salary" row values in row 3
$ws.Range("B3").Value = 1000
$ws.Range("C3").Value = 45905.22928240741

Write-Host "Edit complete"
